# Scheduled market-data refresh for the Maduin Leve profit sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose market snapshot changed.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1725
$ws.Range("I32").Value = 1835
$ws.Range("J32").Value = 1505
$ws.Range("K32").Value = 1835
$ws.Range("L32").Value = 1505
$ws.Range("M32").Value = -1509
$ws.Range("N32").Value = -2157
$ws.Range("H74").Value = 3711
$ws.Range("J74").Value = 2999
$ws.Range("L74").Value = 2999
$ws.Range("N74").Value = -4871
$ws.Range("H76").Value = 2090.3333
$ws.Range("I76").Value = 2987.5
$ws.Range("J76").Value = 296
$ws.Range("K76").Value = 2987.5
$ws.Range("L76").Value = 296
$ws.Range("M76").Value = -2672.5
$ws.Range("N76").Value = -926
$ws.Range("H77").Value = 3711
$ws.Range("J77").Value = 2999
$ws.Range("L77").Value = 14995
$ws.Range("N77").Value = -24355
$ws.Range("H79").Value = 2090.3333
$ws.Range("I79").Value = 2987.5
$ws.Range("J79").Value = 296
$ws.Range("K79").Value = 2987.5
$ws.Range("L79").Value = 296
$ws.Range("M79").Value = -1895.5
$ws.Range("N79").Value = -2480
$ws.Range("H92").Value = 260.7857
$ws.Range("I92").Value = 317.36365
$ws.Range("K92").Value = 317.36365
$ws.Range("M92").Value = 930.63635
$ws.Range("H111").Value = 3762.5
$ws.Range("I111").Value = 3916.6667
$ws.Range("K111").Value = 11750.0001
$ws.Range("M111").Value = -8683.000100000001
$ws.Range("H116").Value = 6103.8696
$ws.Range("I116").Value = 6315.6665
$ws.Range("K116").Value = 6315.6665
$ws.Range("M116").Value = -2873.6665
$ws.Range("H132").Value = 3631.389
$ws.Range("I132").Value = 1905.4166
$ws.Range("J132").Value = 7083.3335
$ws.Range("K132").Value = 5716.2498
$ws.Range("L132").Value = 21250.0005
$ws.Range("M132").Value = -3186.2498
$ws.Range("N132").Value = -26310.0005

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 3100
$ws.Range("J9").Value = 3100
$ws.Range("L9").Value = 3100
$ws.Range("N9").Value = -3440
$ws.Range("H20").Value = 3100
$ws.Range("J20").Value = 3100
$ws.Range("L20").Value = 3100
$ws.Range("N20").Value = -3640
$ws.Range("H32").Value = 5280.5
$ws.Range("I32").Value = 4217.9565
$ws.Range("J32").Value = 17499.75
$ws.Range("K32").Value = 4217.9565
$ws.Range("L32").Value = 17499.75
$ws.Range("M32").Value = -3930.9565
$ws.Range("N32").Value = -18073.75
$ws.Range("H61").Value = 4624.75
$ws.Range("I61").Value = 4624.75
$ws.Range("K61").Value = 4624.75
$ws.Range("M61").Value = -4412.75
$ws.Range("H63").Value = 1500
$ws.Range("I63").Value = 1500
$ws.Range("K63").Value = 1500
$ws.Range("M63").Value = -814
$ws.Range("H66").Value = 1500
$ws.Range("I66").Value = 1500
$ws.Range("K66").Value = 7500
$ws.Range("M66").Value = -4068
$ws.Range("H88").Value = 6176.727
$ws.Range("I88").Value = 5744.5
$ws.Range("K88").Value = 5744.5
$ws.Range("M88").Value = -5338.5
$ws.Range("H91").Value = 6176.727
$ws.Range("I91").Value = 5744.5
$ws.Range("K91").Value = 5744.5
$ws.Range("M91").Value = -4340.5
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H136").Value = 4624.75
$ws.Range("I136").Value = 4624.75
$ws.Range("K136").Value = 13874.25
$ws.Range("M136").Value = -11324.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3709
$ws.Range("J86").Value = 4169
$ws.Range("L86").Value = 4169
$ws.Range("N86").Value = -6415
$ws.Range("H89").Value = 3709
$ws.Range("J89").Value = 4169
$ws.Range("L89").Value = 20845
$ws.Range("N89").Value = -32077
$ws.Range("H105").Value = 3914.8
$ws.Range("I105").Value = 1663.1111
$ws.Range("K105").Value = 1663.1111
$ws.Range("M105").Value = 83.88889999999992
$ws.Range("H134").Value = 4498.4
$ws.Range("I134").Value = 4164.6665
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 12493.9995
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -9958.999500000002
$ws.Range("N134").Value = -20067

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4376
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21880
$ws.Range("N65").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 372.16666
$ws.Range("I40").Value = 378
$ws.Range("K40").Value = 1512
$ws.Range("M40").Value = -1443

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3161.2856
$ws.Range("I80").Value = 2996.5
$ws.Range("K80").Value = 2996.5
$ws.Range("M80").Value = -1998.5
$ws.Range("H83").Value = 3161.2856
$ws.Range("I83").Value = 2996.5
$ws.Range("K83").Value = 14982.5
$ws.Range("M83").Value = -9990.5
$ws.Range("H102").Value = 788.5
$ws.Range("I102").Value = 769.2727
$ws.Range("K102").Value = 769.2727
$ws.Range("M102").Value = 852.7273

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 15855.429
$ws.Range("I40").Value = 24994.5
$ws.Range("J40").Value = 12199.8
$ws.Range("K40").Value = 24994.5
$ws.Range("L40").Value = 12199.8
$ws.Range("M40").Value = -24858.5
$ws.Range("N40").Value = -12471.8

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1120.1
$ws.Range("J81").Value = 900.5
$ws.Range("L81").Value = 1801
$ws.Range("N81").Value = -3923
$ws.Range("H84").Value = 1120.1
$ws.Range("J84").Value = 900.5
$ws.Range("L84").Value = 9005
$ws.Range("N84").Value = -19613
$ws.Range("H107").Value = 763.1429000000001
$ws.Range("I107").Value = 688.4
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 2065.2
$ws.Range("L107").Value = 2850
$ws.Range("M107").Value = -145.1999999999998
$ws.Range("N107").Value = -6690
